$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# New StatQuery text replacing the old stat-count Cypher query used in column C (StatQuery)
$newStatQuery = @"
CALL{
        MATCH (p:participant)-->(s:study)
        OPTIONAL MATCH (samp:sample)-->(p)
        OPTIONAL MATCH (samp)<--(f:file)
        OPTIONAL MATCH (f)<--(g:genomic_info)
        OPTIONAL MATCH (p)<--(diag:diagnosis)
        WITH s, p, samp, f, g, diag
        WHERE f.file_type IN ["JSON"]
        RETURN 
            count(distinct p) AS num_participants
    }
    WITH num_participants
    CALL {
        MATCH (samp:sample)-->(p:participant)-->(s)
        OPTIONAL MATCH (samp)<--(f:file)
        OPTIONAL MATCH (p)<--(diag:diagnosis)
        OPTIONAL MATCH (f)<--(g:genomic_info)
        OPTIONAL MATCH (p)<--(diag:diagnosis)
        WITH s, p, samp, f, g, diag
        WHERE f.file_type IN ["JSON"]
        RETURN 
            count(distinct samp) AS num_samples
    }
    WITH num_participants, num_samples
    CALL {
        MATCH (f:file)-->(s:study)
        OPTIONAL MATCH (f)<--(g:genomic_info)
        OPTIONAL MATCH (samp:sample)<--(f)
        OPTIONAL MATCH (p:participant)<--(samp)
        OPTIONAL MATCH (p)<--(diag:diagnosis)
        WITH s, p, samp, f, g, diag
        WHERE f.file_type IN ["JSON"]
        RETURN 
            count(distinct s) AS num_studies,
            count(distinct f) AS num_files
    }
    RETURN 
        num_studies AS Studies,
        num_participants AS Participants,
        num_samples AS Samples,
        num_files AS ``Files``
"@

# Update the StatQuery column (C) for the three data rows with the new query text
$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Rows grew taller because of the longer query text; match the saved row heights (capped at Excel's max 409.5)
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5

# Update view/selection to match the saved state
$ws.Range("C5").Select()
